$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.082.74'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '2.398.53'
$ws.Range('E3').Value = '  -3.10%  '
$ws.Range('E4').Value = '  +0.06%  '
$origStyle_D5 = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.49'
$ws.Range('D5').Style = $origStyle_D5
$ws.Range('E5').Value = '  -1.43%  '
$origStyle_D6 = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '89.37'
$ws.Range('D6').Style = $origStyle_D6
$ws.Range('E6').Value = '  -3.95%  '
$ws.Range('E7').Value = '  -3.39%  '
$ws.Range('E8').Value = '  +0.12%  '
$origStyle_D9 = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.499'
$ws.Range('D9').Style = $origStyle_D9
$ws.Range('E9').Value = '  -3.52%  '
$origStyle_D10 = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0843'
$ws.Range('D10').Style = $origStyle_D10
$ws.Range('E10').Value = '  -3.17%  '
$origStyle_D11 = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '31.51'
$ws.Range('D11').Style = $origStyle_D11
$ws.Range('E11').Value = '  -5.57%  '
$ws.Range('E12').Value = '  -1.64%  '
$ws.Range('D13').Value = '2.767.60'
$ws.Range('E13').Value = '  -3.07%  '
$origStyle_D14 = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.66'
$ws.Range('D14').Style = $origStyle_D14
$ws.Range('E14').Value = '  -3.85%  '
$origStyle_D15 = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.34'
$ws.Range('D15').Style = $origStyle_D15
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('D16').Value = '2.372.87'
$ws.Range('E16').Value = '  -3.37%  '
$origStyle_D17 = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.775'
$ws.Range('D17').Style = $origStyle_D17
$ws.Range('E17').Value = '  -2.75%  '
$ws.Range('D18').Value = '41.034.48'
$ws.Range('E18').Value = '  -1.62%  '
$ws.Range('D19').Value = '0.0₃0920'
$ws.Range('E19').Value = '  -3.19%  '
$ws.Range('E20').Value = '  -3.82%  '
$origStyle_D21 = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.72'
$ws.Range('D21').Style = $origStyle_D21
$ws.Range('E21').Value = '  -1.97%  '
$ws.Range('E22').Value = '  -3.37%  '
$ws.Range('E23').Value = '  -2.46%  '
$ws.Range('E24').Value = '  -3.05%  '
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('E26').Value = '  -5.44%  '
$origStyle_D27 = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.15'
$ws.Range('D27').Style = $origStyle_D27
$ws.Range('E27').Value = '  -2.46%  '
$ws.Range('E28').Value = '  -1.45%  '
$origStyle_D29 = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.45'
$ws.Range('D29').Style = $origStyle_D29
$ws.Range('E29').Value = '  -3.45%  '
$origStyle_D30 = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.55'
$ws.Range('D30').Style = $origStyle_D30
$ws.Range('E30').Value = '  -5.62%  '
$origStyle_D31 = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '154.48'
$ws.Range('D31').Style = $origStyle_D31
$ws.Range('E31').Value = '  -2.17%  '
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('E33').Value = '  -4.43%  '
$origStyle_D34 = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0741'
$ws.Range('D34').Style = $origStyle_D34
$ws.Range('E34').Value = '  -3.22%  '
$ws.Range('E35').Value = '  -4.81%  '
$ws.Range('E36').Value = '  -1.88%  '
$origStyle_D37 = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.84'
$ws.Range('D37').Style = $origStyle_D37
$ws.Range('E37').Value = '  -2.83%  '
$origStyle_D38 = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.25'
$ws.Range('D38').Style = $origStyle_D38
$ws.Range('E38').Value = '  -7.13%  '
$ws.Range('E39').Value = '  -3.47%  '
$ws.Range('E40').Value = '  -6.26%  '
$ws.Range('E41').Value = '  -2.98%  '
$origStyle_D42 = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.35'
$ws.Range('D42').Style = $origStyle_D42
$ws.Range('E42').Value = '  -8.17%  '
$ws.Range('D43').Value = '1.978.94'
$ws.Range('E43').Value = '  -1.03%  '
$origStyle_D44 = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0273'
$ws.Range('D44').Style = $origStyle_D44
$ws.Range('E44').Value = '  -4.19%  '
$ws.Range('E45').Value = '  -5.72%  '
$origStyle_D46 = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.69'
$ws.Range('D46').Style = $origStyle_D46
$ws.Range('E46').Value = '  +1.30%  '
$ws.Range('E47').Value = '  -6.41%  '
$ws.Range('D48').Value = '2.626.27'
$ws.Range('E48').Value = '  -3.15%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$origStyle_D49 = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '94.53'
$ws.Range('D49').Style = $origStyle_D49
$ws.Range('E49').Value = '  -3.83%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$origStyle_D50 = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.47'
$ws.Range('D50').Style = $origStyle_D50
$ws.Range('E50').Value = '  -1.94%  '
$origStyle_D51 = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.99'
$ws.Range('D51').Style = $origStyle_D51
$ws.Range('E51').Value = '  -1.10%  '
